$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F header (reuse the existing bold/bordered header style from E1) ---
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Mensaje"

# --- Row 2: status flips from "Desconocido" to "Enviado"; new column F stays blank ---
$ws.Range("C2").Value = "Enviado"
$ws.Range("F2").Value = ""

# --- Row 3: clear out the old "Gabriel" contact, keep only the bare phone number ---
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "3339563030"
$ws.Range("C3").Value = ""
$ws.Range("F3").Value = ""

# --- Row 4: brand new test-client row with an auto-sent appointment message ---
$ws.Range("A4").Value = "Prueba Cliente"
$ws.Range("B4").Value = "522205511054"
$ws.Range("C4").Value = "Enviado"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "Cita Agendada"
$ws.Range("F4").Value = "Hola Prueba Cliente, confirmamos tu cita para mañana a las 10:00. 🗓️ ¿Podrías confirmar con un 'SÍ'?"
